$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Column D "error calculation" value swaps on the RM block ---
$ws.Range("D2").Value = -13.5
$ws.Range("D3").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("D11").Value = -15.5
$ws.Range("D13").Value = ""
$ws.Range("D21").Value = -14.3
$ws.Range("D25").Value = ""

# --- 2. Remove the "RM 232" (row 26) and "SC 92" (row 28) data rows entirely ---
# Delete from the bottom up so row numbers of rows still to be removed stay valid.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# --- 3. A few more column C/D "error calculation" corrections on the resulting SC block ---
$ws.Range("C29").Value = ""
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1
